$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the blank separator row (old row 41); this shifts the old
# "Totais" row (old row 42) up to become the new row 41, and updates
# the used range/dimension automatically from A1:C42 to A1:C41.
$ws.Rows.Item(41).Delete()

# Fill in the Qtd_vendida (B) and Valor Total (C) values for each product row.
$ws.Cells.Item(2, 2).Value = 0
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(3, 2).Value = 1
$ws.Cells.Item(3, 3).Value = 25.99
$ws.Cells.Item(4, 2).Value = 2
$ws.Cells.Item(4, 3).Value = 29
$ws.Cells.Item(5, 2).Value = 3
$ws.Cells.Item(5, 3).Value = 129.8593862352
$ws.Cells.Item(6, 2).Value = 11
$ws.Cells.Item(6, 3).Value = 48.45810645244801
$ws.Cells.Item(7, 2).Value = 21
$ws.Cells.Item(7, 3).Value = 125.233750265932
$ws.Cells.Item(8, 2).Value = 0
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(9, 2).Value = 5
$ws.Cells.Item(9, 3).Value = 124.95
$ws.Cells.Item(10, 2).Value = 1
$ws.Cells.Item(10, 3).Value = 24.99
$ws.Cells.Item(11, 2).Value = 3
$ws.Cells.Item(11, 3).Value = 76.5
$ws.Cells.Item(12, 2).Value = 0
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(13, 2).Value = 1
$ws.Cells.Item(13, 3).Value = 25.99
$ws.Cells.Item(14, 2).Value = 1
$ws.Cells.Item(14, 3).Value = 29.99
$ws.Cells.Item(15, 2).Value = 1
$ws.Cells.Item(15, 3).Value = 29.99
$ws.Cells.Item(16, 2).Value = 2
$ws.Cells.Item(16, 3).Value = 9.780000000000001
$ws.Cells.Item(17, 2).Value = 2
$ws.Cells.Item(17, 3).Value = 18.7
$ws.Cells.Item(18, 2).Value = 2
$ws.Cells.Item(18, 3).Value = 29.2
$ws.Cells.Item(19, 2).Value = 6
$ws.Cells.Item(19, 3).Value = 19.74
$ws.Cells.Item(20, 2).Value = 6
$ws.Cells.Item(20, 3).Value = 37.14868489945886
$ws.Cells.Item(21, 2).Value = 6
$ws.Cells.Item(21, 3).Value = 59.15173392000001
$ws.Cells.Item(22, 2).Value = 13
$ws.Cells.Item(22, 3).Value = 94.01340000000002
$ws.Cells.Item(23, 2).Value = 0
$ws.Cells.Item(23, 3).Value = 0
$ws.Cells.Item(24, 2).Value = 2
$ws.Cells.Item(24, 3).Value = 29.1516
$ws.Cells.Item(25, 2).Value = 2
$ws.Cells.Item(25, 3).Value = 51.98
$ws.Cells.Item(26, 2).Value = 2
$ws.Cells.Item(26, 3).Value = 60.58
$ws.Cells.Item(27, 2).Value = 3
$ws.Cells.Item(27, 3).Value = 14.07
$ws.Cells.Item(28, 2).Value = 5
$ws.Cells.Item(28, 3).Value = 46.76598000000001
$ws.Cells.Item(29, 2).Value = 10
$ws.Cells.Item(29, 3).Value = 32.364721584
$ws.Cells.Item(30, 2).Value = 10
$ws.Cells.Item(30, 3).Value = 39.90000000000001
$ws.Cells.Item(31, 2).Value = 15
$ws.Cells.Item(31, 3).Value = 104.85
$ws.Cells.Item(32, 2).Value = 1
$ws.Cells.Item(32, 3).Value = 24.99
$ws.Cells.Item(33, 2).Value = 1
$ws.Cells.Item(33, 3).Value = 30
$ws.Cells.Item(34, 2).Value = 2
$ws.Cells.Item(34, 3).Value = 28
$ws.Cells.Item(35, 2).Value = 3
$ws.Cells.Item(35, 3).Value = 26.97
$ws.Cells.Item(36, 2).Value = 5
$ws.Cells.Item(36, 3).Value = 22.45
$ws.Cells.Item(37, 2).Value = 8
$ws.Cells.Item(37, 3).Value = 68.72000000000001
$ws.Cells.Item(38, 2).Value = 12
$ws.Cells.Item(38, 3).Value = 31.08
$ws.Cells.Item(39, 2).Value = 12
$ws.Cells.Item(39, 3).Value = 47.88
$ws.Cells.Item(40, 2).Value = 12
$ws.Cells.Item(40, 3).Value = 79.08000000000001

# Fill in the totals row (new row 41).
$ws.Cells.Item(41, 2).Value = 192
$ws.Cells.Item(41, 3).Value = 1677.517363357039

Write-Output "done"
